$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Firmware version " paragraph -> append firmware build "SAAFKS00-011-R01E0"
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found = $r1.Find.Execute("Firmware version", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $r1.Paragraphs(1)
    $paraRange = $para.Range
    $insertRange = $d.Range($paraRange.End - 1, $paraRange.End - 1)
    $insertRange.InsertAfter("SAAFKS00-011-R01E0")
}

# ---------------------------------------------------------------------------
# 2. First "2025 NOV LG release" -> "2025 December LG release"
#    (bullet right below "Firmware version ...")
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("NOV", $true, $false, $false, $false, $false, $true, 1, $false, "December", 1)

# ---------------------------------------------------------------------------
# 3. First "SAAFKS00-011-N02" (right after "Android APP 2.0.5.238" hyperlink
#    block) -> "SAAFKS00-011-R01E0.DAT"
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Android APP 2.0.5.238", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scope = $d.Range($anchor.End, $d.Content.End)
$scope.Find.Execute("SAAFKS00-011-N02", $true, $false, $false, $false, $false, $true, 1, $false, "SAAFKS00-011-R01E0.DAT", 1)

# ---------------------------------------------------------------------------
# 4. Second "Scanning framework 2025 NOV LG release" -> "... December ..."
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("NOV", $true, $false, $false, $false, $false, $true, 1, $false, "December", 1)

# ---------------------------------------------------------------------------
# 6. Second "SAAFKS00-011-N02" (right after "Firmware Configuration" heading,
#    "Use the FW SAAFKS00-011-N02   to verify data wedge support")
#    -> "SAAFKS00-011-R01E0" and triple space collapses to double space
# ---------------------------------------------------------------------------
$anchor2 = $d.Content
$anchor2.Find.Execute("Firmware Configuration", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scope2 = $d.Range($anchor2.End, $d.Content.End)
$scope2.Find.Execute("SAAFKS00-011-N02   to verify data wedge support", $true, $false, $false, $false, $false, $true, 1, $false, "SAAFKS00-011-R01E0  to verify data wedge support", 1)

# ---------------------------------------------------------------------------
# 7. Third occurrence, "SAAFKS" + "00-011-N02" (split across two runs) in
#    "Update the latest shared Firmware SAAFKS00-011-N02" -> "SAAFKS00-011-R01E0"
# ---------------------------------------------------------------------------
$r7 = $d.Content
$r7.Find.Execute("SAAFKS00-011-N02", $true, $false, $false, $false, $false, $true, 1, $false, "SAAFKS00-011-R01E0", 1)
